# TC07_Bento_Filter_Chemo-Other.xlsx - "Fixed Bento 80 Test scripts"
# Appends an `order by ... LIMIT 100` clause to the Cypher queries that
# back the CasesTab / SamplesTab / FilesTab rows on the "startup" sheet,
# and nudges the sheet view (scroll position + active cell) to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

# --- B2 (CasesTab query): append an ORDER BY / LIMIT clause -----------
$b2 = $ws.Range("B2").Value2
$ws.Range("B2").Value2 = $b2 + $nl + " order By ss.study_subject_id ASC LIMIT 100"

# --- B3 (SamplesTab query): append an ORDER BY / LIMIT clause ---------
$b3 = $ws.Range("B3").Value2
$ws.Range("B3").Value2 = $b3 + $nl + " order By samp.sample_id ASC LIMIT 100"

# --- B4 (FilesTab query): replace the trailing lower-case "order by"
#     line with a capitalized ORDER BY / LIMIT clause ------------------
$b4 = $ws.Range("B4").Value2
$oldTail = $nl + "    order by f.file_name"
$newTail = $nl + "     order By f.file_name ASC LIMIT 100"
$ws.Range("B4").Value2 = $b4.Replace($oldTail, $newTail)

# --- Row heights grow to fit the extra wrapped line of text -----------
$ws.Rows.Item(2).RowHeight = 331.2
$ws.Rows.Item(3).RowHeight = 360
$ws.Rows.Item(4).RowHeight = 409.6

# --- View state: scroll up one row and move the active cell -----------
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
$ws.Range("B4").Select()
